$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('J2').Value = 7199
$ws.Range('J3').Value = 7607
$ws.Range('J4').Value = 1659
$ws.Range('J5').Value = 594
$ws.Range('J6').Value = 10345
$ws.Range('J7').Value = 27404

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('J2').Value = 72
$ws.Range('J3').Value = 50
$ws.Range('J6').Value = 270
$ws.Range('J7').Value = 409

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('J3').Value = 505
$ws.Range('J6').Value = 635
$ws.Range('J7').Value = 1727

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('J2').Value = 161
$ws.Range('J7').Value = 547

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('J3').Value = 413
$ws.Range('J6').Value = 444
$ws.Range('J7').Value = 1243

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range('J3').Value = 141
$ws.Range('J7').Value = 392

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('J3').Value = 284
$ws.Range('J7').Value = 838

$ws = $wb.Worksheets.Item('New City')
$ws.Range('J3').Value = 188
$ws.Range('J4').Value = 29
$ws.Range('J6').Value = 257
$ws.Range('J7').Value = 689

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('J6').Value = 112
$ws.Range('J7').Value = 421

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range('J2').Value = 34
$ws.Range('J7').Value = 95

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('J4').Value = 126
$ws.Range('J7').Value = 781
$ws.Range('J8').Value = 1727
$ws.Range('J10').Value = 200
$ws.Range('J11').Value = 492
$ws.Range('J14').Value = 147
$ws.Range('J15').Value = 340
$ws.Range('J16').Value = 109
$ws.Range('J17').Value = 36
$ws.Range('J18').Value = 223
$ws.Range('J19').Value = 787
$ws.Range('J20').Value = 581
$ws.Range('J29').Value = 1464
$ws.Range('J30').Value = 95
$ws.Range('J33').Value = 1243
$ws.Range('J37').Value = 838
$ws.Range('J42').Value = 1176
$ws.Range('J43').Value = 230
$ws.Range('J44').Value = 216
$ws.Range('J53').Value = 409
$ws.Range('J54').Value = 543
$ws.Range('J55').Value = 436
$ws.Range('J56').Value = 40
$ws.Range('J58').Value = 16
$ws.Range('J59').Value = 30
$ws.Range('J63').Value = 90
$ws.Range('J64').Value = 180
$ws.Range('J65').Value = 689
$ws.Range('J67').Value = 1020
$ws.Range('J73').Value = 264
$ws.Range('J77').Value = 191
$ws.Range('J79').Value = 753
$ws.Range('J83').Value = 547
$ws.Range('J84').Value = 225
$ws.Range('J85').Value = 1124
$ws.Range('J89').Value = 339
$ws.Range('J93').Value = 118
$ws.Range('J95').Value = 392
$ws.Range('J96').Value = 300
$ws.Range('J99').Value = 421
$ws.Range('J101').Value = 27404

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('J2').Value = 260
$ws.Range('J3').Value = 384
$ws.Range('J7').Value = 1020

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range('J6').Value = 75
$ws.Range('J7').Value = 225

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('J2').Value = 136
$ws.Range('J7').Value = 543

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('J2').Value = 444
$ws.Range('J3').Value = 518
$ws.Range('J5').Value = 55
$ws.Range('J7').Value = 1464

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('J2').Value = 192
$ws.Range('J3').Value = 225
$ws.Range('J6').Value = 303
$ws.Range('J7').Value = 787

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range('J2').Value = 66
$ws.Range('J4').Value = 13
$ws.Range('J6').Value = 86
$ws.Range('J7').Value = 216

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range('J3').Value = 28
$ws.Range('J7').Value = 147

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range('J5').Value = 9
$ws.Range('J6').Value = 75

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('J2').Value = 243
$ws.Range('J3').Value = 236
$ws.Range('J4').Value = 51
$ws.Range('J6').Value = 625
$ws.Range('J7').Value = 1176

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range('J6').Value = 114
$ws.Range('J7').Value = 200

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('J6').Value = 246
$ws.Range('J7').Value = 436

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('J6').Value = 108
$ws.Range('J7').Value = 300

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('J3').Value = 253
$ws.Range('J7').Value = 753

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range('J6').Value = 63
$ws.Range('J7').Value = 180

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('J2').Value = 162
$ws.Range('J3').Value = 194
$ws.Range('J6').Value = 168
$ws.Range('J7').Value = 581

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range('J3').Value = 47
$ws.Range('J7').Value = 223

$ws = $wb.Worksheets.Item('Burnside')
$ws.Range('J6').Value = 6
$ws.Range('J7').Value = 36

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range('J2').Value = 35
$ws.Range('J4').Value = 10
$ws.Range('J7').Value = 118

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('J6').Value = 249
$ws.Range('J7').Value = 781

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range('J2').Value = 35
$ws.Range('J6').Value = 49

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range('J6').Value = 157
$ws.Range('J7').Value = 340

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('J6').Value = 233
$ws.Range('J7').Value = 492

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range('J3').Value = 65
$ws.Range('J7').Value = 264

$ws = $wb.Worksheets.Item('Montclare')
$ws.Range('J2').Value = 15
$ws.Range('J7').Value = 30

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('J6').Value = 104
$ws.Range('J7').Value = 339

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range('J6').Value = 136
$ws.Range('J7').Value = 230

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('J2').Value = 302
$ws.Range('J6').Value = 318
$ws.Range('J7').Value = 1124

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range('J2').Value = 71
$ws.Range('J7').Value = 191

$ws = $wb.Worksheets.Item('Magnificent Mile')
$ws.Range('J3').Value = 10
$ws.Range('J6').Value = 19
$ws.Range('J7').Value = 40

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range('J5').Value = 3
$ws.Range('J7').Value = 126

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range('J6').Value = 85
$ws.Range('J7').Value = 109

$ws = $wb.Worksheets.Item('Millenium Park')
$ws.Range('J2').Value = 3
$ws.Range('J7').Value = 16
